$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart")

# --- Insert two new rows at row 12 (rows 12/13), pushing everything down ---
$ws.Range("A12:A13").EntireRow.Insert()

# Copy formatting for the two new rows from rows 4 and 5 (same visual pattern)
$ws.Range("A4:R4").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 43.5

$ws.Range("A5:R5").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = 29

Write-Output "step1 done"

# --- Fill new row 12 content ---
$ws.Range("B12").Value2 = 7
$ws.Range("C12").Value2 = "As a player, I need to log in to the site in order to play the game"
$ws.Range("D12").Value2 = "After log in page, bring up game board"
$ws.Range("E12").Value2 = "Patrick Garcia"
$ws.Range("F12").Value2 = 5
$ws.Range("H12").Value2 = 6

# --- Fill new row 13 content ---
$ws.Range("D13").Value2 = "Connect site to game logic"
$ws.Range("E13").Value2 = "Patrick Garcia"
$ws.Range("F13").Value2 = 5
$ws.Range("H13").Value2 = 6

Write-Output "step2 done"

# --- Update rows 4-11: assignee, estimate adjustments, and actuals ---
$ws.Range("E4").Value2 = "Mark Pratt"
$ws.Range("G4").Value2 = 2
$ws.Range("H4:P4").ClearContents()

$ws.Range("E5").Value2 = "Mark Pratt"
$ws.Range("G5").Value2 = 1
$ws.Range("H5:P5").ClearContents()

$ws.Range("E6").Value2 = "Patrick Garcia"
$ws.Range("F6").Value2 = 3
$ws.Range("G6").Value2 = 3
$ws.Range("H6:P6").ClearContents()

$ws.Range("E7").Value2 = "Patrick Garcia"
$ws.Range("G7").Value2 = 2
$ws.Range("H7:P7").ClearContents()

$ws.Range("E8").Value2 = "Patrick Garcia"
$ws.Range("G8:P8").ClearContents()

$ws.Range("E9").Value2 = "Mark Pratt"
$ws.Range("G9:P9").ClearContents()

$ws.Range("E10").Value2 = "Mark Pratt"
$ws.Range("G10:P10").ClearContents()

$ws.Range("E11").Value2 = "Mark Pratt"
$ws.Range("G11:P11").ClearContents()

Write-Output "step3 done"

# --- Update totals rows (now row14 = Ideal/Hours trend, row15 = Actual) ---
$ws.Range("F14").Formula = '=SUM(F4:F13)'
$ws.Range("G14").Formula = '=F14-$F$14/10'
$ws.Range("H14:O14").Formula = '=G14-$F$14/10'
$ws.Range("P14").Formula = '=IF((O14-$F$14/10) >= 1, (O14-$F$14/10), 0)'

$ws.Range("F15").Formula = '=SUM(F4:F13)'
$ws.Range("G15").Formula = '=F15 - SUM(G4:G13)'
$ws.Range("H15").Formula = '=G15 - SUM(H4:H13)'
$ws.Range("I15:P15").Formula = '=H15 - SUM(I4:I13)'

Write-Output "step4 done"

# --- Update title text for Sprint 2 ---
$ws.Range("A1").Value2 = "Project Title: Minesweeper CLC Milestone`r`nRelease #:`r`nSprint #: 2"

Write-Output "step5 done"
